# Weekly update: insert a new price record for "Inferno" (Ají) at the top
# of the recent-dates block (row 13), pushing the existing rows 13-16 down
# to 14-17 — mirrors the source data feed's newest-first ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13:16 shift down to 14:17
# and inherit their formatting (date style, etc.) automatically.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44553
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112021
$ws.Range("G13").Value = "Ají"
$ws.Range("H13").Value = "Inferno"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 35
$ws.Range("K13").Value = 45000
$ws.Range("L13").Value = 45000
$ws.Range("M13").Value = 45000
$ws.Range("N13").Value = "`$/caja 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1800
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
